$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1797
$ws.Range("F3").Value = 10478
$ws.Range("F4").Value = 15
$ws.Range("F10").Value = 23
$ws.Range("F16").Value = 138
$ws.Range("F17").Value = 492
$ws.Range("F19").Value = 398
$ws.Range("F20").Value = 117
$ws.Range("F21").Value = 1210
$ws.Range("F22").Value = 1146
$ws.Range("F23").Value = 1304
$ws.Range("F24").Value = 242
$ws.Range("F25").Value = 1476
$ws.Range("F26").Value = 437
$ws.Range("F27").Value = 746
$ws.Range("F28").Value = 278
$ws.Range("F31").Value = 894
$ws.Range("F32").Value = 281
$ws.Range("F33").Value = 763
$ws.Range("F35").Value = 871
$ws.Range("F36").Value = 160161
$ws.Range("F37").Value = 856
$ws.Range("F38").Value = 538
$ws.Range("F40").Value = 857
$ws.Range("F41").Value = 787
$ws.Range("F42").Value = 1469
$ws.Range("F44").Value = 745

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 130
$ws.Range("F14").Value = 1257
$ws.Range("F16").Value = 2368
$ws.Range("F18").Value = 363
$ws.Range("F20").Value = 146
$ws.Range("F22").Value = 49
$ws.Range("F32").Value = 211
$ws.Range("F42").Value = 100
$ws.Range("F44").Value = 1

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 847
$ws.Range("F6").Value = 2613
$ws.Range("F7").Value = 4344
$ws.Range("F8").Value = 95
$ws.Range("F10").Value = 458
$ws.Range("F11").Value = 470
$ws.Range("F12").Value = 345
$ws.Range("F13").Value = 347

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1797
$ws.Range("F3").Value = 847
$ws.Range("F5").Value = 4344
$ws.Range("F6").Value = 95
$ws.Range("F8").Value = 470
$ws.Range("F10").Value = 23
$ws.Range("F12").Value = 347
$ws.Range("F13").Value = 347
$ws.Range("F15").Value = 130
$ws.Range("F19").Value = 1257
$ws.Range("F20").Value = 492
$ws.Range("F22").Value = 398
$ws.Range("F23").Value = 117
$ws.Range("F24").Value = 2368
$ws.Range("F26").Value = 1146
$ws.Range("F27").Value = 1304
$ws.Range("F28").Value = 146
$ws.Range("F29").Value = 49
$ws.Range("F30").Value = 1476
$ws.Range("F31").Value = 746
$ws.Range("F32").Value = 278
$ws.Range("F34").Value = 894
$ws.Range("F35").Value = 763
$ws.Range("F37").Value = 871
$ws.Range("F39").Value = 856
$ws.Range("F40").Value = 538
$ws.Range("F41").Value = 857
$ws.Range("F42").Value = 787
$ws.Range("F44").Value = 1469
$ws.Range("F48").Value = 745
$ws.Range("F50").Value = 100
